{"js": "// \"download articles with pandoc title blocks\"\n//\n// Replace the old byline block:\n//   P1: italic \"On Pilgrimage\" + \",\" + \" \" + <br/> + \" \" + \"September ==================\"\n//   P2: bold \"By Dorothy Day\"\n// with a pandoc-style title block:\n//   P1: paragraph styled \"Title\", plain text \"September\"\n//   P2: plain text \"% Dorothy Day\"\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0];\nconst secondPara = paragraphs.items[1];\n\n// First paragraph: wipe all runs/formatting, replace with plain \"September\",\n// then apply the built-in \"Title\" paragraph style.\nconst firstRange = firstPara.getRange();\nfirstRange.clear();\nfirstRange.insertText(\"September\", Word.InsertLocation.replace);\nfirstPara.style = \"Title\";\n\n// Second paragraph: wipe the bold run, replace with plain \"% Dorothy Day\".\nconst secondRange = secondPara.getRange();\nsecondRange.clear();\nsecondRange.insertText(\"% Dorothy Day\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# \"download articles with pandoc title blocks\"\n#\n# Replace the old byline block:\n#   P1: italic \"On Pilgrimage\" + \",\" + \" \" + <br/> + \" \" + \"September ==================\"\n#   P2: bold \"By Dorothy Day\"\n# with a pandoc-style title block:\n#   P1: paragraph styled \"Title\", plain text \"September\"\n#   P2: plain text \"% Dorothy Day\"\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: strip all runs/formatting, then insert plain \"September\"\n#     and apply the built-in \"Title\" paragraph style. ---\n$p1Range = $d.Paragraphs(1).Range\n$p1Body = $d.Range($p1Range.Start, $p1Range.End - 1)\n$p1Body.Delete()\n$d.Paragraphs(1).Range.InsertAfter(\"September\")\n$d.Paragraphs(1).Style = \"Title\"\n\n# --- Paragraph 2: strip the bold run, then insert plain \"% Dorothy Day\". ---\n$p2Range = $d.Paragraphs(2).Range\n$p2Body = $d.Range($p2Range.Start, $p2Range.End - 1)\n$p2Body.Delete()\n$d.Paragraphs(2).Range.InsertAfter(\"% Dorothy Day\")\n"}
